# Sentralitetsindeks_SSB.xlsx — refresh Sentralitetsindeks data (Power BI
# re-pull) and re-sort the table by Kommunenummer instead of
# Sentralitetsklasse.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Sentralitetsindeks" (col D), "Intervall" (col E) and
# "Sentralitetsklasse" (col F) values, in the CURRENT (pre-sort) row order
# of the sheet (rows 2..42, still ordered/grouped by Sentralitetsklasse as
# they were before this edit). The row identity (Kommunenummer in col A)
# does not change — only these three derived/measured columns do.
$D_VALUES = @(793,686,686,713,679,670,613,653,603,631,650,605,584,426,505,400,373,480,559,494,453,429,494,567,362,323,338,488,478,448,547,403,537,528,366,370,522,497,517,442,405)
$E_VALUES = @("775-869","670-774","670-774","670-774","670-774","670-774","565-669","565-669","565-669","565-669","565-669","565-669","565-669","0-564","0-564","0-564","0-564","0-564","0-564","0-564","0-564","0-564","0-564","565-669","0-564","0-564","0-564","0-564","0-564","0-564","0-564","0-564","0-564","0-564","0-564","0-564","0-564","0-564","0-564","0-564","0-564")
$F_VALUES = @(3,4,4,4,4,4,5,5,5,5,5,5,5,6,6,6,6,6,6,6,6,6,6,5,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6)

for ($i = 0; $i -lt 41; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 4).Value2 = $D_VALUES[$i]
    $ws.Cells.Item($r, 5).Value2 = $E_VALUES[$i]
    $ws.Cells.Item($r, 6).Value2 = $F_VALUES[$i]
}

# Re-sort the table by Kommunenummer (column A) ascending instead of by
# Sentralitetsklasse (column F).
$lo = $ws.ListObjects.Item(1)
$lo.Sort.SortFields.Clear()
$lo.Sort.SortFields.Add($ws.Range("A2:A42"))
$lo.Sort.Header = 1
$lo.Sort.Apply()

# Match the saved selection from the edited workbook.
$ws.Range("H13").Select()

Write-Output "done"
